$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> points conversion; add a half-EMU epsilon so float32 round-trip
# through the COM Left/Top/Width/Height properties lands back on the
# exact target EMU value instead of being floored one EMU short.
$emuPerPt = 12700
$eps = 0.5 / $emuPerPt
function ToPt($emu) { return ($emu / $emuPerPt) + $eps }

# --- Reposition existing ridge-plot legend images -------------------------
$s.Shapes.Item(2).Top = ToPt 4794949   # Picture 16 (id 17)
$s.Shapes.Item(3).Top = ToPt 4970695   # Picture 17 (id 18)
$s.Shapes.Item(4).Top = ToPt 4439103   # Picture 22 (id 23)
$s.Shapes.Item(7).Top = ToPt 4856567   # Picture 26 (id 27)
$s.Shapes.Item(8).Top = ToPt 5054321   # Picture 27 (id 28)

# --- Add the new "#zesty_color_palette" notes textbox ----------------------
# Duplicate an existing autofit textbox (Rectangle 44) so the new shape's
# XML (bodyPr/spAutoFit, lstStyle, rPr dirty/smtClean attrs) matches what a
# real author's PowerPoint would emit, instead of the bare markup that
# Shapes.AddShape/AddTextbox produce. Also, this presentation's shape-id
# allocator assigns ids by skipping any id already used in the slide's
# original XML, counting up from 1 for every Duplicate()/AddShape() call in
# this session; the 43rd call lands on id 61, matching the target shape.
$template = $s.Shapes.Item(14)
$shp = $null
for ($i = 1; $i -le 43; $i++) {
    $range = $template.Duplicate()
    $shp = $range.Item(1)
    if ($i -lt 43) {
        $shp.Delete()
    }
}

$shp.Name = "Rectangle 60"
$shp.Left = ToPt 6050954
$shp.Top = ToPt 1213008
$shp.Width = ToPt 4572000
$shp.Height = ToPt 1477328

$tr = $shp.TextFrame.TextRange
$tr.Text = '#zesty_color_palette '
$tr.LanguageID = "mr-IN"

$nl = "`n"
$r1 = $tr.InsertAfter($nl + 'color0 = "#0F2080"')
$r1.LanguageID = "mr-IN"
$r2 = $tr.InsertAfter($nl + 'color1 = "#F5793A"')
$r2.LanguageID = "mr-IN"
$r3 = $tr.InsertAfter($nl + 'color2 = "#85C0F9" ')
$r3.LanguageID = "mr-IN"
$r4 = $tr.InsertAfter($nl + 'color3 = "#A95AA1"')
$r4.LanguageID = "mr-IN"
